$d = $word.ActiveDocument

# These character styles have both <w:color> and <w:b/> (bold only) in their
# rPr. The color/bold order needs to be fixed to "b before color" to match
# wml.xsd (CT_RPr sequence). Re-assigning Font.Bold forces the engine to
# regenerate <w:rPr> with schema-correct element ordering.
$boldOnlyStyles = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
foreach ($styleName in $boldOnlyStyles) {
    $s = $d.Styles.Item($styleName)
    $s.Font.Bold = $True
}

# These character styles have both <w:color> and <w:i/> (italic only) in
# their rPr. Re-assigning Font.Italic forces the engine to regenerate
# <w:rPr> with schema-correct element ordering (i before color).
$italicOnlyStyles = @("CommentTok", "DocumentationTok")
foreach ($styleName in $italicOnlyStyles) {
    $s = $d.Styles.Item($styleName)
    $s.Font.Italic = $True
}

# These character styles have <w:color>, <w:b/> and <w:i/> in their rPr.
# Re-assigning both Font.Bold and Font.Italic forces the engine to
# regenerate <w:rPr> with schema-correct element ordering (b, i, then color).
$boldItalicStyles = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")
foreach ($styleName in $boldItalicStyles) {
    $s = $d.Styles.Item($styleName)
    $s.Font.Bold = $True
    $s.Font.Italic = $True
}
